# Update the NATMI ligand-receptor edge table (Calr -> Lrp1) with
# recomputed TPM-based values. The underlying sending/target cluster
# average & total expression values changed, which cascades into the
# derived specificity (I,J,O,P) and edge weight/specificity (Q,R,S,T)
# columns for every row. Values below are the new, recomputed figures
# for each (Sending cluster, Target cluster) row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 63.18739966666666
$ws.Cells.Item(2, 8).Value = 189.562199
$ws.Cells.Item(2, 9).Value = 0.09596345243430386
$ws.Cells.Item(2, 10).Value = 0.09988075390087989
$ws.Cells.Item(2, 13).Value = 2.906846333333333
$ws.Cells.Item(2, 14).Value = 8.720538999999999
$ws.Cells.Item(2, 15).Value = 0.005520525738044089
$ws.Cells.Item(2, 16).Value = 0.005624540846623205
$ws.Cells.Item(2, 17).Value = 183.6760610339178
$ws.Cells.Item(2, 18).Value = 1653.084549305261
$ws.Cells.Item(2, 19).Value = 0.0005297687090751441
$ws.Cells.Item(2, 20).Value = 0.000561783380107019

# Row 3
$ws.Cells.Item(3, 7).Value = 63.18739966666666
$ws.Cells.Item(3, 8).Value = 189.562199
$ws.Cells.Item(3, 9).Value = 0.09596345243430386
$ws.Cells.Item(3, 10).Value = 0.09988075390087989
$ws.Cells.Item(3, 15).Value = 0.3528665483720876
$ws.Cells.Item(3, 16).Value = 0.3595150912979765
$ws.Cells.Item(3, 17).Value = 11740.39226535382
$ws.Cells.Item(3, 18).Value = 105663.5303881844
$ws.Cells.Item(3, 19).Value = 0.03386229223036182
$ws.Cells.Item(3, 20).Value = 0.03590863835758556

# Row 4
$ws.Cells.Item(4, 7).Value = 63.18739966666666
$ws.Cells.Item(4, 8).Value = 189.562199
$ws.Cells.Item(4, 9).Value = 0.09596345243430386
$ws.Cells.Item(4, 10).Value = 0.09988075390087989
$ws.Cells.Item(4, 13).Value = 137.0717086666666
$ws.Cells.Item(4, 14).Value = 411.2151259999999
$ws.Cells.Item(4, 15).Value = 0.2603191943704447
$ws.Cells.Item(4, 16).Value = 0.2652240042658267
$ws.Cells.Item(4, 17).Value = 8661.204838513562
$ws.Cells.Item(4, 18).Value = 77950.84354662205
$ws.Cells.Item(4, 19).Value = 0.02498112862670447
$ws.Cells.Item(4, 20).Value = 0.02649077349868096

# Row 5
$ws.Cells.Item(5, 7).Value = 63.18739966666666
$ws.Cells.Item(5, 8).Value = 189.562199
$ws.Cells.Item(5, 9).Value = 0.09596345243430386
$ws.Cells.Item(5, 10).Value = 0.09988075390087989
$ws.Cells.Item(5, 13).Value = 29.2127365
$ws.Cells.Item(5, 14).Value = 58.425473
$ws.Cells.Item(5, 15).Value = 0.05547925319534149
$ws.Cells.Item(5, 16).Value = 0.03768304451958546
$ws.Cells.Item(5, 17).Value = 1845.876856582521
$ws.Cells.Item(5, 18).Value = 11075.26113949513
$ws.Cells.Item(5, 19).Value = 0.005323980675101853
$ws.Cells.Item(5, 20).Value = 0.003763810895896615

# Row 6
$ws.Cells.Item(6, 7).Value = 63.18739966666666
$ws.Cells.Item(6, 8).Value = 189.562199
$ws.Cells.Item(6, 9).Value = 0.09596345243430386
$ws.Cells.Item(6, 10).Value = 0.09988075390087989
$ws.Cells.Item(6, 13).Value = 171.5584106666666
$ws.Cells.Item(6, 14).Value = 514.6752319999999
$ws.Cells.Item(6, 15).Value = 0.3258144783240821
$ws.Cells.Item(6, 16).Value = 0.331953319069988
$ws.Cells.Item(6, 17).Value = 10840.3298609728
$ws.Cells.Item(6, 18).Value = 97562.96874875516
$ws.Cells.Item(6, 19).Value = 0.03126628219306058
$ws.Cells.Item(6, 20).Value = 0.03315574776860973

# Row 7
$ws.Cells.Item(7, 9).Value = 0.3063997713314046
$ws.Cells.Item(7, 10).Value = 0.3189072441572365
$ws.Cells.Item(7, 13).Value = 2.906846333333333
$ws.Cells.Item(7, 14).Value = 8.720538999999999
$ws.Cells.Item(7, 15).Value = 0.005520525738044089
$ws.Cells.Item(7, 16).Value = 0.005624540846623205
$ws.Cells.Item(7, 17).Value = 586.4555898337798
$ws.Cells.Item(7, 18).Value = 5278.100308504018
$ws.Cells.Item(7, 19).Value = 0.001691487823765842
$ws.Cells.Item(7, 20).Value = 0.001793706821046417

# Row 8
$ws.Cells.Item(8, 9).Value = 0.3063997713314046
$ws.Cells.Item(8, 10).Value = 0.3189072441572365
$ws.Cells.Item(8, 15).Value = 0.3528665483720876
$ws.Cells.Item(8, 16).Value = 0.3595150912979765
$ws.Cells.Item(8, 19).Value = 0.1081182297317097
$ws.Cells.Item(8, 20).Value = 0.114651966998775

# Row 9
$ws.Cells.Item(9, 9).Value = 0.3063997713314046
$ws.Cells.Item(9, 10).Value = 0.3189072441572365
$ws.Cells.Item(9, 13).Value = 137.0717086666666
$ws.Cells.Item(9, 14).Value = 411.2151259999999
$ws.Cells.Item(9, 15).Value = 0.2603191943704447
$ws.Cells.Item(9, 16).Value = 0.2652240042658267
$ws.Cells.Item(9, 17).Value = 27654.18619960327
$ws.Cells.Item(9, 18).Value = 248887.6757964294
$ws.Cells.Item(9, 19).Value = 0.07976174162827972
$ws.Cells.Item(9, 20).Value = 0.08458185628476195

# Row 10
$ws.Cells.Item(10, 9).Value = 0.3063997713314046
$ws.Cells.Item(10, 10).Value = 0.3189072441572365
$ws.Cells.Item(10, 13).Value = 29.2127365
$ws.Cells.Item(10, 14).Value = 58.425473
$ws.Cells.Item(10, 15).Value = 0.05547925319534149
$ws.Cells.Item(10, 16).Value = 0.03768304451958546
$ws.Cells.Item(10, 17).Value = 5893.662977058972
$ws.Cells.Item(10, 18).Value = 35361.97786235383
$ws.Cells.Item(10, 19).Value = 0.01699883049268973
$ws.Cells.Item(10, 20).Value = 0.01201739587919545

# Row 11
$ws.Cells.Item(11, 9).Value = 0.3063997713314046
$ws.Cells.Item(11, 10).Value = 0.3189072441572365
$ws.Cells.Item(11, 13).Value = 171.5584106666666
$ws.Cells.Item(11, 14).Value = 514.6752319999999
$ws.Cells.Item(11, 15).Value = 0.3258144783240821
$ws.Cells.Item(11, 16).Value = 0.331953319069988
$ws.Cells.Item(11, 17).Value = 34611.87052261305
$ws.Cells.Item(11, 18).Value = 311506.8347035174
$ws.Cells.Item(11, 19).Value = 0.09982948165495964
$ws.Cells.Item(11, 20).Value = 0.1058623181734577

# Row 12
$ws.Cells.Item(12, 7).Value = 170.2928416666667
$ws.Cells.Item(12, 8).Value = 510.878525
$ws.Cells.Item(12, 9).Value = 0.2586257560429799
$ws.Cells.Item(12, 10).Value = 0.2691830570543736
$ws.Cells.Item(12, 13).Value = 2.906846333333333
$ws.Cells.Item(12, 14).Value = 8.720538999999999
$ws.Cells.Item(12, 15).Value = 0.005520525738044089
$ws.Cells.Item(12, 16).Value = 0.005624540846623205
$ws.Cells.Item(12, 17).Value = 495.0151223916638
$ws.Cells.Item(12, 18).Value = 4455.136101524975
$ws.Cells.Item(12, 19).Value = 0.001427750142756382
$ws.Cells.Item(12, 20).Value = 0.001514031099621229

# Row 13
$ws.Cells.Item(13, 7).Value = 170.2928416666667
$ws.Cells.Item(13, 8).Value = 510.878525
$ws.Cells.Item(13, 9).Value = 0.2586257560429799
$ws.Cells.Item(13, 10).Value = 0.2691830570543736
$ws.Cells.Item(13, 15).Value = 0.3528665483720876
$ws.Cells.Item(13, 16).Value = 0.3595150912979765
$ws.Cells.Item(13, 17).Value = 31640.87732198848
$ws.Cells.Item(13, 18).Value = 284767.8958978963
$ws.Cells.Item(13, 19).Value = 0.09126037785500792
$ws.Cells.Item(13, 20).Value = 0.09677537133277155

# Row 14
$ws.Cells.Item(14, 7).Value = 170.2928416666667
$ws.Cells.Item(14, 8).Value = 510.878525
$ws.Cells.Item(14, 9).Value = 0.2586257560429799
$ws.Cells.Item(14, 10).Value = 0.2691830570543736
$ws.Cells.Item(14, 13).Value = 137.0717086666666
$ws.Cells.Item(14, 14).Value = 411.2151259999999
$ws.Cells.Item(14, 15).Value = 0.2603191943704447
$ws.Cells.Item(14, 16).Value = 0.2652240042658267
$ws.Cells.Item(14, 17).Value = 23342.33078095212
$ws.Cells.Item(14, 18).Value = 210080.9770285691
$ws.Cells.Item(14, 19).Value = 0.0673252484565557
$ws.Cells.Item(14, 20).Value = 0.07139380827247745

# Row 15
$ws.Cells.Item(15, 7).Value = 170.2928416666667
$ws.Cells.Item(15, 8).Value = 510.878525
$ws.Cells.Item(15, 9).Value = 0.2586257560429799
$ws.Cells.Item(15, 10).Value = 0.2691830570543736
$ws.Cells.Item(15, 13).Value = 29.2127365
$ws.Cells.Item(15, 14).Value = 58.425473
$ws.Cells.Item(15, 15).Value = 0.05547925319534149
$ws.Cells.Item(15, 16).Value = 0.03768304451958546
$ws.Cells.Item(15, 17).Value = 4974.719911444554
$ws.Cells.Item(15, 18).Value = 29848.31946866732
$ws.Cells.Item(15, 19).Value = 0.0143483638023451
$ws.Cells.Item(15, 20).Value = 0.01014363712289807

# Row 16
$ws.Cells.Item(16, 7).Value = 170.2928416666667
$ws.Cells.Item(16, 8).Value = 510.878525
$ws.Cells.Item(16, 9).Value = 0.2586257560429799
$ws.Cells.Item(16, 10).Value = 0.2691830570543736
$ws.Cells.Item(16, 13).Value = 171.5584106666666
$ws.Cells.Item(16, 14).Value = 514.6752319999999
$ws.Cells.Item(16, 15).Value = 0.3258144783240821
$ws.Cells.Item(16, 16).Value = 0.331953319069988
$ws.Cells.Item(16, 17).Value = 29215.16926424364
$ws.Cells.Item(16, 18).Value = 262936.5233781928
$ws.Cells.Item(16, 19).Value = 0.08426401578631483
$ws.Cells.Item(16, 20).Value = 0.08935620922660527

# Row 17
$ws.Cells.Item(17, 7).Value = 77.473122
$ws.Cells.Item(17, 8).Value = 154.946244
$ws.Cells.Item(17, 9).Value = 0.1176593481802354
$ws.Cells.Item(17, 10).Value = 0.08164152846121862
$ws.Cells.Item(17, 13).Value = 2.906846333333333
$ws.Cells.Item(17, 14).Value = 8.720538999999999
$ws.Cells.Item(17, 15).Value = 0.005520525738044089
$ws.Cells.Item(17, 16).Value = 0.005624540846623205
$ws.Cells.Item(17, 17).Value = 225.202460617586
$ws.Cells.Item(17, 18).Value = 1351.214763705516
$ws.Cells.Item(17, 19).Value = 0.0006495414599504802
$ws.Cells.Item(17, 20).Value = 0.0004591961116108752

# Row 18
$ws.Cells.Item(18, 7).Value = 77.473122
$ws.Cells.Item(18, 8).Value = 154.946244
$ws.Cells.Item(18, 9).Value = 0.1176593481802354
$ws.Cells.Item(18, 10).Value = 0.08164152846121862
$ws.Cells.Item(18, 15).Value = 0.3528665483720876
$ws.Cells.Item(18, 16).Value = 0.3595150912979765
$ws.Cells.Item(18, 17).Value = 14394.71867967114
$ws.Cells.Item(18, 18).Value = 86368.31207802685
$ws.Cells.Item(18, 19).Value = 0.04151804807606932
$ws.Cells.Item(18, 20).Value = 0.02935136155844136

# Row 19
$ws.Cells.Item(19, 7).Value = 77.473122
$ws.Cells.Item(19, 8).Value = 154.946244
$ws.Cells.Item(19, 9).Value = 0.1176593481802354
$ws.Cells.Item(19, 10).Value = 0.08164152846121862
$ws.Cells.Item(19, 13).Value = 137.0717086666666
$ws.Cells.Item(19, 14).Value = 411.2151259999999
$ws.Cells.Item(19, 15).Value = 0.2603191943704447
$ws.Cells.Item(19, 16).Value = 0.2652240042658267
$ws.Cells.Item(19, 17).Value = 10619.37320828112
$ws.Cells.Item(19, 18).Value = 63716.23924968674
$ws.Cells.Item(19, 19).Value = 0.03062898672843051
$ws.Cells.Item(19, 20).Value = 0.02165329309286686

# Row 20
$ws.Cells.Item(20, 7).Value = 77.473122
$ws.Cells.Item(20, 8).Value = 154.946244
$ws.Cells.Item(20, 9).Value = 0.1176593481802354
$ws.Cells.Item(20, 10).Value = 0.08164152846121862
$ws.Cells.Item(20, 13).Value = 29.2127365
$ws.Cells.Item(20, 14).Value = 58.425473
$ws.Cells.Item(20, 15).Value = 0.05547925319534149
$ws.Cells.Item(20, 16).Value = 0.03768304451958546
$ws.Cells.Item(20, 17).Value = 2263.201898818353
$ws.Cells.Item(20, 18).Value = 9052.807595273413
$ws.Cells.Item(20, 19).Value = 0.006527652768490119
$ws.Cells.Item(20, 20).Value = 0.003076501351651105

# Row 21
$ws.Cells.Item(21, 7).Value = 77.473122
$ws.Cells.Item(21, 8).Value = 154.946244
$ws.Cells.Item(21, 9).Value = 0.1176593481802354
$ws.Cells.Item(21, 10).Value = 0.08164152846121862
$ws.Cells.Item(21, 13).Value = 171.5584106666666
$ws.Cells.Item(21, 14).Value = 514.6752319999999
$ws.Cells.Item(21, 15).Value = 0.3258144783240821
$ws.Cells.Item(21, 16).Value = 0.331953319069988
$ws.Cells.Item(21, 17).Value = 13291.16567970477
$ws.Cells.Item(21, 18).Value = 79746.9940782286
$ws.Cells.Item(21, 19).Value = 0.03833511914729492
$ws.Cells.Item(21, 20).Value = 0.02710117634664842

# Row 22
$ws.Cells.Item(22, 7).Value = 145.7496183333334
$ws.Cells.Item(22, 8).Value = 437.248855
$ws.Cells.Item(22, 9).Value = 0.2213516720110761
$ws.Cells.Item(22, 10).Value = 0.2303874164262914
$ws.Cells.Item(22, 13).Value = 2.906846333333333
$ws.Cells.Item(22, 14).Value = 8.720538999999999
$ws.Cells.Item(22, 15).Value = 0.005520525738044089
$ws.Cells.Item(22, 16).Value = 0.005624540846623205
$ws.Cells.Item(22, 17).Value = 423.6717436369828
$ws.Cells.Item(22, 18).Value = 3813.045692732845
$ws.Cells.Item(22, 19).Value = 0.001221977602496239
$ws.Cells.Item(22, 20).Value = 0.001295823434237666

# Row 23
$ws.Cells.Item(23, 7).Value = 145.7496183333334
$ws.Cells.Item(23, 8).Value = 437.248855
$ws.Cells.Item(23, 9).Value = 0.2213516720110761
$ws.Cells.Item(23, 10).Value = 0.2303874164262914
$ws.Cells.Item(23, 15).Value = 0.3528665483720876
$ws.Cells.Item(23, 16).Value = 0.3595150912979765
$ws.Cells.Item(23, 17).Value = 27080.67907186925
$ws.Cells.Item(23, 18).Value = 243726.1116468232
$ws.Cells.Item(23, 19).Value = 0.07810760047893886
$ws.Cells.Item(23, 20).Value = 0.08282775305040309

# Row 24
$ws.Cells.Item(24, 7).Value = 145.7496183333334
$ws.Cells.Item(24, 8).Value = 437.248855
$ws.Cells.Item(24, 9).Value = 0.2213516720110761
$ws.Cells.Item(24, 10).Value = 0.2303874164262914
$ws.Cells.Item(24, 13).Value = 137.0717086666666
$ws.Cells.Item(24, 14).Value = 411.2151259999999
$ws.Cells.Item(24, 15).Value = 0.2603191943704447
$ws.Cells.Item(24, 16).Value = 0.2652240042658267
$ws.Cells.Item(24, 17).Value = 19978.14922246453
$ws.Cells.Item(24, 18).Value = 179803.3430021807
$ws.Cells.Item(24, 19).Value = 0.05762208893047423
$ws.Cells.Item(24, 20).Value = 0.0611042731170395

# Row 25
$ws.Cells.Item(25, 7).Value = 145.7496183333334
$ws.Cells.Item(25, 8).Value = 437.248855
$ws.Cells.Item(25, 9).Value = 0.2213516720110761
$ws.Cells.Item(25, 10).Value = 0.2303874164262914
$ws.Cells.Item(25, 13).Value = 29.2127365
$ws.Cells.Item(25, 14).Value = 58.425473
$ws.Cells.Item(25, 15).Value = 0.05547925319534149
$ws.Cells.Item(25, 16).Value = 0.03768304451958546
$ws.Cells.Item(25, 17).Value = 4257.745195347236
$ws.Cells.Item(25, 18).Value = 25546.47117208342
$ws.Cells.Item(25, 19).Value = 0.01228042545671467
$ws.Cells.Item(25, 20).Value = 0.008681699269944212

# Row 26
$ws.Cells.Item(26, 7).Value = 145.7496183333334
$ws.Cells.Item(26, 8).Value = 437.248855
$ws.Cells.Item(26, 9).Value = 0.2213516720110761
$ws.Cells.Item(26, 10).Value = 0.2303874164262914
$ws.Cells.Item(26, 13).Value = 171.5584106666666
$ws.Cells.Item(26, 14).Value = 514.6752319999999
$ws.Cells.Item(26, 15).Value = 0.3258144783240821
$ws.Cells.Item(26, 16).Value = 0.331953319069988
$ws.Cells.Item(26, 17).Value = 25004.57287653993
$ws.Cells.Item(26, 18).Value = 225041.1558888594
$ws.Cells.Item(26, 19).Value = 0.07211957954245209
$ws.Cells.Item(26, 20).Value = 0.07647786755466691
